$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.022792
$ws.Range("H2").Value = 0.06837600000000001
$ws.Range("I2").Value = 0.001916327914826657
$ws.Range("J2").Value = 0.001916327914826657
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.023286
$ws.Range("N2").Value = 0.069858
$ws.Range("O2").Value = 0.009310710475795457
$ws.Range("P2").Value = 0.009310710475795458
$ws.Range("Q2").Value = 0.0005307345120000001
$ws.Range("R2").Value = 0.004776610608000001
$ws.Range("S2").Value = 0.00001784237439163582
$ws.Range("T2").Value = 0.00001784237439163582
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.022792
$ws.Range("H3").Value = 0.06837600000000001
$ws.Range("I3").Value = 0.001916327914826657
$ws.Range("J3").Value = 0.001916327914826657
$ws.Range("O3").Value = 0.05314667307834813
$ws.Range("P3").Value = 0.05314667307834814
$ws.Range("Q3").Value = 0.003029497445333333
$ws.Range("R3").Value = 0.027265477008
$ws.Range("S3").Value = 0.0001018464532002049
$ws.Range("T3").Value = 0.0001018464532002049
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.022792
$ws.Range("H4").Value = 0.06837600000000001
$ws.Range("I4").Value = 0.001916327914826657
$ws.Range("J4").Value = 0.001916327914826657
$ws.Range("M4").Value = 2.344785333333334
$ws.Range("N4").Value = 7.034356000000001
$ws.Range("O4").Value = 0.9375426164458565
$ws.Range("P4").Value = 0.9375426164458565
$ws.Range("Q4").Value = 0.05344234731733335
$ws.Range("R4").Value = 0.4809811258560001
$ws.Range("S4").Value = 0.001796639087234816
$ws.Range("T4").Value = 0.001796639087234816
$ws.Range("I5").Value = 0.3701235913233977
$ws.Range("J5").Value = 0.3701235913233977
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.023286
$ws.Range("N5").Value = 0.069858
$ws.Range("O5").Value = 0.009310710475795457
$ws.Range("P5").Value = 0.009310710475795458
$ws.Range("Q5").Value = 0.102507176408
$ws.Range("R5").Value = 0.922564587672
$ws.Range("S5").Value = 0.003446113599073796
$ws.Range("T5").Value = 0.003446113599073796
$ws.Range("I6").Value = 0.3701235913233977
$ws.Range("J6").Value = 0.3701235913233977
$ws.Range("O6").Value = 0.05314667307834813
$ws.Range("P6").Value = 0.05314667307834814
$ws.Range("S6").Value = 0.01967083750664875
$ws.Range("T6").Value = 0.01967083750664875
$ws.Range("I7").Value = 0.3701235913233977
$ws.Range("J7").Value = 0.3701235913233977
$ws.Range("M7").Value = 2.344785333333334
$ws.Range("N7").Value = 7.034356000000001
$ws.Range("O7").Value = 0.9375426164458565
$ws.Range("P7").Value = 0.9375426164458565
$ws.Range("Q7").Value = 10.32196701034489
$ws.Range("R7").Value = 92.89770309310401
$ws.Range("S7").Value = 0.3470066402176752
$ws.Range("T7").Value = 0.3470066402176752
$ws.Range("G8").Value = 7.468693666666667
$ws.Range("H8").Value = 22.406081
$ws.Range("I8").Value = 0.6279600807617757
$ws.Range("J8").Value = 0.6279600807617757
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.023286
$ws.Range("N8").Value = 0.069858
$ws.Range("O8").Value = 0.009310710475795457
$ws.Range("P8").Value = 0.009310710475795458
$ws.Range("Q8").Value = 0.173916000722
$ws.Range("R8").Value = 1.565244006498
$ws.Range("S8").Value = 0.005846754502330026
$ws.Range("T8").Value = 0.005846754502330027
$ws.Range("G9").Value = 7.468693666666667
$ws.Range("H9").Value = 22.406081
$ws.Range("I9").Value = 0.6279600807617757
$ws.Range("J9").Value = 0.6279600807617757
$ws.Range("O9").Value = 0.05314667307834813
$ws.Range("P9").Value = 0.05314667307834814
$ws.Range("Q9").Value = 0.992733783044222
$ws.Range("R9").Value = 8.934604047397999
$ws.Range("S9").Value = 0.03337398911849918
$ws.Range("T9").Value = 0.03337398911849919
$ws.Range("G10").Value = 7.468693666666667
$ws.Range("H10").Value = 22.406081
$ws.Range("I10").Value = 0.6279600807617757
$ws.Range("J10").Value = 0.6279600807617757
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 2.344785333333334
$ws.Range("N10").Value = 7.034356000000001
$ws.Range("O10").Value = 0.9375426164458565
$ws.Range("P10").Value = 0.9375426164458565
$ws.Range("Q10").Value = 17.51248336875956
$ws.Range("R10").Value = 157.612350318836
$ws.Range("S10").Value = 0.5887393371409465
$ws.Range("T10").Value = 0.5887393371409465
